$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Label")
$ws2.Activate()
try {
  $excel.Goto($ws2.Range("C29"), $true)
  Write-Host "Goto worked"
} catch {
  Write-Host "Goto failed: $_"
}
